# Weekly data refresh: insert the newest price record at the top of the
# data block (row 40), pushing all existing historical rows down by one.
# The former last row (63) survives as the new last row (64), and the
# sheet's used range grows from R63 to R64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 40..63 down to 41..64, inheriting formatting from the row
# above the insertion point (keeps the date-formatted style on column D).
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with this week's record.
$ws.Cells.Item(40, 1).Value  = 11
$ws.Cells.Item(40, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(40, 3).Value  = "Bíobío"
$ws.Cells.Item(40, 4).Value  = 44582
$ws.Cells.Item(40, 5).Value  = 8
$ws.Cells.Item(40, 6).Value  = 100112001
$ws.Cells.Item(40, 7).Value  = "Berenjena"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 100
$ws.Cells.Item(40, 11).Value = 10000
$ws.Cells.Item(40, 12).Value = 11000
$ws.Cells.Item(40, 13).Value = 10500
$ws.Cells.Item(40, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(40, 15).Value = "Región Metropolitana"
$ws.Cells.Item(40, 16).Value = 175
$ws.Cells.Item(40, 17).Value = 60
$ws.Cells.Item(40, 18).Value = "Hortaliza"
